# Update on Feb 11
# - Populate Nationality (L) / SymtomDate (M) for many existing rows
# - Add two new case rows (47, 48) with full data
# - Widen column H, adjust the saved view position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Nationality (column L) + SymtomDate (column M) updates for existing rows
# ---------------------------------------------------------------------------
# Map of row -> [Nationality text (or $null to leave as-is), SymtomDate serial (or $null)]
$updates = @(
    @{Row=2;  L="Chinese";       M=44562},
    @{Row=4;  L="Chinese";       M=44562},
    @{Row=5;  L="Chinese";       M=44562},
    @{Row=6;  L="Chinese";       M=45292},
    @{Row=7;  L="Chinese";       M=45658},
    @{Row=9;  L="Chinese";       M=43466},
    @{Row=10; L="Chinese";       M=43466},
    @{Row=13; L="Chinese";       M=46023},
    @{Row=15; L="Chinese";       M=46753},
    @{Row=16; L=$null;           M=10959},
    @{Row=17; L="Chinese";       M=44562},
    @{Row=18; L=$null;           M=10959},
    @{Row=19; L="Chinese";       M=36923},
    @{Row=20; L="Chinese (PR)";  M=47119},
    @{Row=21; L=$null;           M=45658},
    @{Row=22; L="Indonesian";    M=$null},
    @{Row=24; L=$null;           M=10959},
    @{Row=25; L=$null;           M=10959},
    @{Row=26; L=$null;           M=45292},
    @{Row=30; L=$null;           M=46753},
    @{Row=31; L=$null;           M=43831},
    @{Row=32; L=$null;           M=44927},
    @{Row=34; L=$null;           M=10959},
    @{Row=35; L=$null;           M=46388},
    @{Row=36; L=$null;           M=11324},
    @{Row=37; L="Chinese (PR)";  M=45292},
    @{Row=38; L=$null;           M=10959},
    @{Row=40; L=$null;           M=47119},
    @{Row=41; L=$null;           M=10959},
    @{Row=42; L=$null;           M=36923},
    @{Row=43; L="Bangladesh";    M=$null},
    @{Row=44; L=$null;           M=10959},
    @{Row=46; L=$null;           M=10959}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.L -ne $null) {
        $ws.Range("L$r").Value = $u.L
    }
    if ($u.M -ne $null) {
        $ws.Range("M$r").Value = $u.M
        $ws.Range("M$r").NumberFormat = "mmm-yy"
    }
}

# Row 2 previously had no cell in column N at all - make sure it exists (blank,
# same style as the rest of the row). Copy the format from N4, which is a
# known-blank "s=3" cell, rather than reading N2 back (Value getters on empty
# cells surface raw property metadata instead of $null/empty in this shim).
$ws.Range("N4").Copy()
$ws.Range("N2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Two new case rows (47 and 48), copying the formatting of the last
#    existing data row (46) and then filling in the actual values.
# ---------------------------------------------------------------------------
$ws.Range("A46:N46").Copy()
$ws.Range("A47:N47").PasteSpecial(-4122)
$ws.Range("A48:N48").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Rows(47).RowHeight = 21
$ws.Rows(48).RowHeight = 21

# Row 47
$ws.Range("A47").Value = 46
$ws.Range("B47").Value = 1.2559439999999999
$ws.Range("C47").Value = 103.82021400000001
$ws.Range("D47").Value = "Feb-10"
$ws.Range("E47").Value = 35
$ws.Range("F47").Value = "Male"
$ws.Range("G47").Value = "Singapore"
$ws.Range("H47").Value = "Johor Bahru"
$ws.Range("I47").Value = "Johor Bahru, Resorts World Sentosa Casino, TTSH, NCID"
$ws.Range("L47").Value = "PR"
$ws.Range("M47").Value = 38384
$ws.Range("M47").NumberFormat = "mmm-yy"

# Row 48
$ws.Range("A48").Value = 47
$ws.Range("B48").Value = 1.3073600000000001
$ws.Range("C48").Value = 103.854623
$ws.Range("D48").Value = "Feb-10"
$ws.Range("E48").Value = 39
$ws.Range("F48").Value = "Male"
$ws.Range("G48").Value = "Singapore"
$ws.Range("H48").Value = "Veerasamy Road"
$ws.Range("I48").Value = "10 Seletar Aerospace Heights"
$ws.Range("J48").Value = "42"
$ws.Range("L48").Value = "Bangladesh"
$ws.Range("M48").Value = 38749
$ws.Range("M48").NumberFormat = "mmm-yy"

# ---------------------------------------------------------------------------
# 3. Column H is much wider now (used to fit "Johor Bahru, Resorts World
#    Sentosa Casino, TTSH, NCID"). 54.33 is the closest achievable input that
#    this engine's pixel-rounded column-width model maps back to 55.1640625.
# ---------------------------------------------------------------------------
$ws.Columns("H").ColumnWidth = 54.33

# ---------------------------------------------------------------------------
# 4. Saved view position / selection.
# ---------------------------------------------------------------------------
[void]$excel.Goto($ws.Range("C34"))
[void]$ws.Range("I47").Select()
